# Update the ECO report figures (B/D/E columns) for each row.
# Source cells are stored as text (e.g. "6.00", "-19.00", "24.00%"), and a
# plain Range.Value assignment of a numeric-looking string would silently be
# coerced to a real number by Excel (losing the fixed ".00"/"%" formatting
# and flipping the cell's stored type). To keep the cells as text - exactly
# like the original workbook - each cell is briefly marked with a "@" (Text)
# number format before the value is written, then restored to the "Normal"
# style so no extra formatting/style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, [string]$Address, [string]$NewValue)

    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

# Row 2 - Beatrice Chege
Set-TextValue $ws "B2" "5.00"
Set-TextValue $ws "D2" "-20.00"
Set-TextValue $ws "E2" "20.00%"

# Row 3 - Glayds  Bundotich
Set-TextValue $ws "B3" "1.00"
Set-TextValue $ws "D3" "-19.00"
Set-TextValue $ws "E3" "5.00%"

# Row 4 - Jane Gichohi
Set-TextValue $ws "B4" "4.00"
Set-TextValue $ws "D4" "-21.00"
Set-TextValue $ws "E4" "16.00%"

# Row 5 - Mirriam Makau
Set-TextValue $ws "B5" "2.00"
Set-TextValue $ws "D5" "-18.00"
Set-TextValue $ws "E5" "10.00%"

# Row 6 - Victor Njogu
Set-TextValue $ws "B6" "2.00"
Set-TextValue $ws "D6" "-23.00"
Set-TextValue $ws "E6" "8.00%"

# Row 7 - KD Totals
Set-TextValue $ws "B7" "14.00"
Set-TextValue $ws "D7" "-101.00"
Set-TextValue $ws "E7" "11.80%"
